$d = $word.ActiveDocument

# --- Edit 1: "If a player plays multiple games ..." -> split into
#     "If a player plays m" | "any" | " games there is not much improvement in "
$rng1 = $d.Content
$rng1.Find.Execute("If a player plays multiple games there is not much improvement in ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start1 = $rng1.Start

# Force run boundaries by dropping temporary bookmarks around the
# substring that needs to change ("multiple" -> "any"), then clean them
# up once the text has been replaced in its own run.
$d.Bookmarks.Add("__TmpSplitA", $d.Range($start1 + 19, $start1 + 19))
$d.Bookmarks.Add("__TmpSplitB", $d.Range($start1 + 26, $start1 + 26))

$midRng1 = $d.Range($start1 + 19, $start1 + 26)
$midRng1.Text = "any"

$d.Bookmarks("__TmpSplitA").Delete()
$d.Bookmarks("__TmpSplitB").Delete()

# --- Edit 2: "To use this system in other games values of constants ..."
#     -> split into "To use this system in other gam" | [_GoBack] | "es values ..."
$rng2 = $d.Content
$rng2.Find.Execute("To use this system in other games values of constants", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start2 = $rng2.Start

# --- Edit 3: drop the old "_GoBack" bookmark before re-adding it at its
#     new location (Word only ever keeps one "_GoBack").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$d.Bookmarks.Add("_GoBack", $d.Range($start2 + 31, $start2 + 31))
